# Insert a new data row at 285 (pushes old rows 285-332 down to 286-333),
# then populate the new row 285 with the new weekly record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 285..332 down by one, inheriting formatting from the row above.
$ws.Rows("285:285").Insert()

# Populate the newly inserted row 285 with the new record.
$ws.Cells.Item(285, 1).Value = 11
$ws.Cells.Item(285, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(285, 3).Value = "Bíobío"
$ws.Cells.Item(285, 4).Value = 45209
$ws.Cells.Item(285, 5).Value = 8
$ws.Cells.Item(285, 6).Value = 100112003
$ws.Cells.Item(285, 7).Value = "Ajo"
$ws.Cells.Item(285, 8).Value = "Chino"
$ws.Cells.Item(285, 9).Value = "Primera"
$ws.Cells.Item(285, 10).Value = 220
$ws.Cells.Item(285, 11).Value = 19000
$ws.Cells.Item(285, 12).Value = 20000
$ws.Cells.Item(285, 13).Value = 19545
$ws.Cells.Item(285, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(285, 15).Value = "China"
$ws.Cells.Item(285, 16).Value = 1954
$ws.Cells.Item(285, 17).Value = 10
$ws.Cells.Item(285, 18).Value = "Hortaliza"
